# Update the marksheet totals for correct/total marks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row - marks awarded per correct answer (Right column)
$ws.Range("B11").Value = 5

# "Total" row - total marks earned for correct answers (Right column)
$ws.Range("B12").Value = 40

# Corr/total marks summary text
$ws.Range("E12").Value = "40/140"
